$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.015.29"
$ws.Range("E2").Value = "  +2.45%  "

$ws.Range("D3").Value = "1.820.13"
$ws.Range("E3").Value = "  +3.14%  "

$ws.Range("D4").Value = "1.014"
$ws.Range("E4").Value = "  +1.23%  "

$ws.Range("D5").Value = "312.31"
$ws.Range("E5").Value = "  +2.37%  "

$ws.Range("E6").Value = "  +0.76%  "

$ws.Range("D7").Value = "0.4313"
$ws.Range("E7").Value = "  +0.73%  "

$ws.Range("D8").Value = "0.3701"
$ws.Range("E8").Value = "  +2.34%  "

$ws.Range("D9").Value = "0.07282"
$ws.Range("E9").Value = "  +3.34%  "

$ws.Range("D10").Value = "2.186.60"
$ws.Range("E10").Value = "  +23.30%  "

$ws.Range("D11").Value = "0.8699"
$ws.Range("E11").Value = "  +4.33%  "

$ws.Range("D12").Value = "21.32"
$ws.Range("E12").Value = "  +5.51%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "5.423"
$ws.Range("E13").Value = "  +3.57%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "6.653"
$ws.Range("E14").Value = "  +3.81%  "

$ws.Range("D15").Value = "0.06974"
$ws.Range("E15").Value = "  +2.69%  "

$ws.Range("D16").Value = "81.20"
$ws.Range("E16").Value = "  +2.52%  "

$ws.Range("E17").Value = "  +1.10%  "

$ws.Range("D18").Value = "0.000008934"
$ws.Range("E18").Value = "  +3.45%  "

$ws.Range("E19").Value = "  +0.64%  "

$ws.Range("D20").Value = "15.34"
$ws.Range("E20").Value = "  +2.38%  "

$ws.Range("D21").Value = "27.050.69"
$ws.Range("E21").Value = "  +2.58%  "

$ws.Range("D22").Value = "5.218"
$ws.Range("E22").Value = "  +4.37%  "

$ws.Range("B23").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D23").Value = "2.389.48"
$ws.Range("E23").Value = "  +20.12%  "

$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "10.99"
$ws.Range("E24").Value = "  -0.97%  "

$ws.Range("D25").Value = "154.78"
$ws.Range("E25").Value = "  +1.48%  "

$ws.Range("D26").Value = "1.893"
$ws.Range("E26").Value = "  +0.79%  "

$ws.Range("D27").Value = "18.39"
$ws.Range("E27").Value = "  +1.53%  "

$ws.Range("D28").Value = "5.248"
$ws.Range("E28").Value = "  +4.20%  "

$ws.Range("D29").Value = "1.929"
$ws.Range("E29").Value = "  +14.35%  "

$ws.Range("D30").Value = "115.11"
$ws.Range("E30").Value = "  +0.54%  "

$ws.Range("D31").Value = "0.08968"
$ws.Range("E31").Value = "  +0.75%  "

$ws.Range("D32").Value = "1.174"
$ws.Range("E32").Value = "  +6.09%  "

$ws.Range("D33").Value = "0.7472"
$ws.Range("E33").Value = "  +3.28%  "

$ws.Range("D34").Value = "4.444"
$ws.Range("E34").Value = "  +2.97%  "

$ws.Range("D35").Value = "2.822"
$ws.Range("E35").Value = "  +2.56%  "

$ws.Range("E36").Value = "  +0.79%  "

$ws.Range("D37").Value = "1.127"
$ws.Range("E37").Value = "  +5.31%  "

$ws.Range("D38").Value = "0.05243"
$ws.Range("E38").Value = "  +2.93%  "

$ws.Range("E39").Value = "  +2.16%  "

$ws.Range("D40").Value = "0.5132"
$ws.Range("E40").Value = "  +4.79%  "

$ws.Range("D41").Value = "2.747"
$ws.Range("E41").Value = "  +10.16%  "

$ws.Range("D42").Value = "0.1653"
$ws.Range("E42").Value = "  +3.50%  "

$ws.Range("D43").Value = "6.507"
$ws.Range("E43").Value = "  +4.65%  "

$ws.Range("D44").Value = "8.354"
$ws.Range("E44").Value = "  +4.29%  "

$ws.Range("D45").Value = "107.45"
$ws.Range("E45").Value = "  +2.43%  "

$ws.Range("E46").Value = "  +3.72%  "

$ws.Range("E47").Value = "  +0.72%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.06316"
$ws.Range("E48").Value = "  +2.10%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "1.653"
$ws.Range("E49").Value = "  +5.22%  "

$ws.Range("D50").Value = "0.4575"
$ws.Range("E50").Value = "  +2.33%  "

$ws.Range("D51").Value = "1.817"
$ws.Range("E51").Value = "  +5.87%  "
